$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap the displayed country names for rows 202 and 203 ---
# (Before: row202="Timor Oriental", row203="Santa Lucia"
#  After:  row202="Santa Lucia",   row203="Timor Oriental")
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Update the "Datos actualizados..." timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 05:46"

# --- Update Kazajistan stats (row 32) ---
$ws.Range("B32").Value = 104718
$ws.Range("C32").Value = 175
$ws.Range("D32").Value = 92598
$ws.Range("E32").Value = 10705

# --- Update Belgica stats (row 40) ---
$ws.Range("B40").Value = 81936
$ws.Range("C40").Value = 468
$ws.Range("D40").Value = 18225
$ws.Range("E40").Value = 53719
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 9992

# --- Update Honduras stats (row 51) ---
$ws.Range("B51").Value = 54511
$ws.Range("C51").Value = 528
$ws.Range("D51").Value = 8532
$ws.Range("E51").Value = 44325
$ws.Range("G51").Value = 11
$ws.Range("H51").Value = 1654

# --- Update Vietnam stats (row 159) ---
$ws.Range("D159").Value = 568
$ws.Range("E159").Value = 421

# --- Update Mongolia stats (row 183) ---
$ws.Range("D183").Value = 289
$ws.Range("E183").Value = 9
